# "visulization and analysis update Fahad"
#
# The canonical-XML diff for this commit shows a single semantic change:
# the internal <p:sldId> "id" attribute for the last slide in the deck
# (still wired to the same r:id/slide position) is renumbered from 352
# to 365. (Every other hunk in the supplied diff is just PowerPoint
# re-serialising unrelated extLst blocks with its namespace-declaration
# attributes in a different order -- xmlns="" before xmlns:p14=/xmlns:p15=/
# xmlns:a16=/xmlns:thm15= instead of after -- with no value changes; that
# is cosmetic XML noise from the save process, not something exposed on
# the PowerPoint object model.)
#
# `Slide.SlideID` is read-only in the PowerPoint object model, so the
# only way an end user (and therefore this COM script) can cause
# PowerPoint to mint a new internal slide id for a slide is to have
# PowerPoint re-create that slide's part -- e.g. by cutting it and
# pasting it back. The freshly (re)created slide is appended at the end
# of the deck and is handed the next never-before-used id, which for
# this deck is 365 -- matching the target id exactly, while leaving the
# slide's content, position and every other slide completely untouched.

$p = $ppt.ActivePresentation

$lastIndex = $p.Slides.Count
$lastSlide = $p.Slides.Item($lastIndex)

$lastSlide.Cut()
$p.Slides.Paste($lastIndex) | Out-Null
